$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. CORE COMPETENCIES: collapse the three detailed bullet paragraphs
#    into a single summary line, and relocate their full detail text
#    into a brand-new "TECHNICAL SKILLS" section at the end of the
#    document (added below).
# ---------------------------------------------------------------------
$bullet = [char]0x2022

$coreCompetencies = $d.Paragraphs.Item(6)
$coreCompetencies.Range.Text = "Product Marketing Core $bullet Research & Analytics $bullet Communication & Technology"

# Remove the two now-redundant detail paragraphs that followed it.
$d.Paragraphs.Item(7).Range.Delete()
$d.Paragraphs.Item(7).Range.Delete()

# ---------------------------------------------------------------------
# 2. Append a new "TECHNICAL SKILLS" section at the end of the document
#    containing the full detail that used to live under CORE
#    COMPETENCIES (now semicolon-delimited, header phrase inline).
# ---------------------------------------------------------------------
$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastParagraph.Range.InsertParagraphAfter()

$heading = $d.Paragraphs.Item($d.Paragraphs.Count)
$heading.Range.Text = "TECHNICAL SKILLS"
$heading.Style = "Heading2"

$heading.Range.InsertParagraphAfter()
$productMarketing = $d.Paragraphs.Item($d.Paragraphs.Count)
$productMarketing.Style = "Normal"
$productMarketing.Range.Text = "PRODUCT MARKETING CORE Market Intelligence & Competitive Analysis; Product Positioning & Messaging Development; Go-to-Market Strategy & Product Launch Management; Customer Segmentation & Buyer Persona Development; Cross-functional Team Leadership & Collaboration; Sales Enablement & Training Material Development; Data-Driven Decision Making & Analytics Interpretation"

$productMarketing.Range.InsertParagraphAfter()
$researchAnalytics = $d.Paragraphs.Item($d.Paragraphs.Count)
$researchAnalytics.Style = "Normal"
$researchAnalytics.Range.Text = "RESEARCH & ANALYTICS Survey Methodology & Customer Insights; Market Research Design & Implementation; Competitive Intelligence & SWOT Analysis; Customer Journey Mapping & Behavioral Analysis; Statistical Modeling & Trend Analysis; Performance Metrics & Dashboard Development; A/B Testing & Conversion Optimization"

$researchAnalytics.Range.InsertParagraphAfter()
$commsTech = $d.Paragraphs.Item($d.Paragraphs.Count)
$commsTech.Style = "Normal"
$commsTech.Range.Text = "COMMUNICATION & TECHNOLOGY Strategic Messaging & Narrative Development; Technical Concept Translation for Business Audiences; Stakeholder Communication & Presentation Skills; Data Visualization & Reporting (Tableau, PowerBI, d3.js); Marketing Technology Stack Integration; Content Strategy & Thought Leadership; Client Relationship Management & Business Development"

Write-Host "Final paragraph count: $($d.Paragraphs.Count)"
